$wb = $excel.ActiveWorkbook

# --- Update selections on the pre-existing sheets (done first so the
# new sheet added later ends up being the active/selected tab) ---

$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("B38").Select()

$ws2 = $wb.Worksheets.Item("W5")
$ws2.Activate()
$ws2.Range("A30:O36").Select()

$ws3 = $wb.Worksheets.Item("W4")
$ws3.Activate()
$ws3.Range("M13").Select()

# --- Add the new "Sheet1" worksheet at the end of the workbook ---

$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Sheet1"

# Column widths (approximate the author's manual column resizing)
$newSheet.Columns.Item(1).ColumnWidth = 25.666666666666668
$newSheet.Columns.Item(2).ColumnWidth = 12.5
$newSheet.Columns.Item(3).ColumnWidth = 5.0
$newSheet.Columns.Item(4).ColumnWidth = 4.833333333333333
$newSheet.Columns.Item(5).ColumnWidth = 4.666666666666667
$newSheet.Columns.Item(6).ColumnWidth = 4.666666666666667
$newSheet.Columns.Item(11).ColumnWidth = 6.333333333333333
$newSheet.Columns.Item(12).ColumnWidth = 7.0

$newSheet.Range("A1").Value = "Device"
$newSheet.Range("B1").HorizontalAlignment = -4131
$newSheet.Range("B1").Value = "Date"
$newSheet.Range("C1").Value = "CO"
$newSheet.Range("D1").Value = "NO2"
$newSheet.Range("E1").Value = "SO2"
$newSheet.Range("F1").Value = "O3"
$newSheet.Range("G1").Value = "rawCO"
$newSheet.Range("H1").Value = "rawNO2"
$newSheet.Range("I1").Value = "rawSO2"
$newSheet.Range("J1").Value = "rawO3"
$newSheet.Range("K1").Value = "AQI"
$newSheet.Range("L1").Value = "Temp"
$newSheet.Range("M1").Value = "Humid"
$newSheet.Range("N1").Value = "Press"
$newSheet.Range("O1").Value = "Place"
$newSheet.Range("A2").Value = "00:A0:50:1A:D6:9E - AirifyW5"
$newSheet.Range("B2").HorizontalAlignment = -4131
$newSheet.Range("B2").NumberFormat = "m/d/yy h:mm"
$newSheet.Range("B2").Value = 43836.245011574072
$newSheet.Range("C2").Value = 42
$newSheet.Range("D2").Value = 83
$newSheet.Range("E2").Value = 137
$newSheet.Range("F2").Value = 33
$newSheet.Range("G2").Value = 1695233
$newSheet.Range("H2").Value = 5627469
$newSheet.Range("I2").Value = 1684771
$newSheet.Range("J2").Value = 5626695
$newSheet.Range("K2").Value = 137
$newSheet.Range("L2").Value = 25.8
$newSheet.Range("M2").Value = 43.66
$newSheet.Range("N2").Value = 989.17
$newSheet.Range("O2").Value = "HOME"
$newSheet.Range("A3").Value = "00:A0:50:1A:D6:9E - AirifyW5"
$newSheet.Range("B3").HorizontalAlignment = -4131
$newSheet.Range("B3").NumberFormat = "m/d/yy h:mm"
$newSheet.Range("B3").Value = 43836.245717592596
$newSheet.Range("C3").Value = 40
$newSheet.Range("D3").Value = 72
$newSheet.Range("E3").Value = 105
$newSheet.Range("F3").Value = 41
$newSheet.Range("G3").Value = 1694385
$newSheet.Range("H3").Value = 5627696
$newSheet.Range("I3").Value = 1681627
$newSheet.Range("J3").Value = 5626107
$newSheet.Range("K3").Value = 105
$newSheet.Range("L3").Value = 25.8
$newSheet.Range("M3").Value = 43.66
$newSheet.Range("N3").Value = 989.17
$newSheet.Range("O3").Value = "HOME"
$newSheet.Range("A4").Value = "00:A0:50:1A:D6:9E - AirifyW5"
$newSheet.Range("B4").HorizontalAlignment = -4131
$newSheet.Range("B4").NumberFormat = "m/d/yy h:mm"
$newSheet.Range("B4").Value = 43836.246423611112
$newSheet.Range("C4").Value = 37
$newSheet.Range("D4").Value = 71
$newSheet.Range("E4").Value = 107
$newSheet.Range("F4").Value = 39
$newSheet.Range("G4").Value = 1693278
$newSheet.Range("H4").Value = 5627709
$newSheet.Range("I4").Value = 1681834
$newSheet.Range("J4").Value = 5626311
$newSheet.Range("K4").Value = 107
$newSheet.Range("L4").Value = 25.8
$newSheet.Range("M4").Value = 43.66
$newSheet.Range("N4").Value = 989.17
$newSheet.Range("O4").Value = "HOME"
$newSheet.Range("A5").Value = "00:A0:50:1A:D6:9E - AirifyW5"
$newSheet.Range("B5").HorizontalAlignment = -4131
$newSheet.Range("B5").NumberFormat = "m/d/yy h:mm"
$newSheet.Range("B5").Value = 43836.247129629628
$newSheet.Range("C5").Value = 36
$newSheet.Range("D5").Value = 79
$newSheet.Range("E5").Value = 106
$newSheet.Range("F5").Value = 37
$newSheet.Range("G5").Value = 1692743
$newSheet.Range("H5").Value = 5627555
$newSheet.Range("I5").Value = 1681669
$newSheet.Range("J5").Value = 5626424
$newSheet.Range("K5").Value = 106
$newSheet.Range("L5").Value = 24.26
$newSheet.Range("M5").Value = 45.25
$newSheet.Range("N5").Value = 983.05
$newSheet.Range("O5").Value = "HOME"
$newSheet.Range("A6").Value = "00:A0:50:1A:D6:9E - AirifyW5"
$newSheet.Range("B6").HorizontalAlignment = -4131
$newSheet.Range("B6").NumberFormat = "m/d/yy h:mm"
$newSheet.Range("B6").Value = 43836.247835648152
$newSheet.Range("C6").Value = 36
$newSheet.Range("D6").Value = 79
$newSheet.Range("E6").Value = 106
$newSheet.Range("F6").Value = 37
$newSheet.Range("G6").Value = 1692743
$newSheet.Range("H6").Value = 5627555
$newSheet.Range("I6").Value = 1681669
$newSheet.Range("J6").Value = 5626424
$newSheet.Range("K6").Value = 106
$newSheet.Range("L6").Value = 24.26
$newSheet.Range("M6").Value = 45.25
$newSheet.Range("N6").Value = 983.05
$newSheet.Range("O6").Value = "HOME"
$newSheet.Range("A7").Value = "00:A0:50:1A:D6:9E - AirifyW5"
$newSheet.Range("B7").HorizontalAlignment = -4131
$newSheet.Range("B7").NumberFormat = "m/d/yy h:mm"
$newSheet.Range("B7").Value = 43836.248541666668
$newSheet.Range("C7").Value = 35
$newSheet.Range("D7").Value = 70
$newSheet.Range("E7").Value = 98
$newSheet.Range("F7").Value = 38
$newSheet.Range("G7").Value = 1692558
$newSheet.Range("H7").Value = 5627730
$newSheet.Range("I7").Value = 1681122
$newSheet.Range("J7").Value = 5626379
$newSheet.Range("K7").Value = 98
$newSheet.Range("L7").Value = 24.26
$newSheet.Range("M7").Value = 45.25
$newSheet.Range("N7").Value = 983.05
$newSheet.Range("O7").Value = "HOME"
$newSheet.Range("A8").Value = "00:A0:50:1A:D6:9E - AirifyW5"
$newSheet.Range("B8").HorizontalAlignment = -4131
$newSheet.Range("B8").NumberFormat = "m/d/yy h:mm"
$newSheet.Range("B8").Value = 43836.249247685184
$newSheet.Range("C8").Value = 35
$newSheet.Range("D8").Value = 74
$newSheet.Range("E8").Value = 88
$newSheet.Range("F8").Value = 41
$newSheet.Range("G8").Value = 1692270
$newSheet.Range("H8").Value = 5627653
$newSheet.Range("I8").Value = 1680764
$newSheet.Range("J8").Value = 5626132
$newSheet.Range("K8").Value = 88
$newSheet.Range("L8").Value = 24.26
$newSheet.Range("M8").Value = 45.25
$newSheet.Range("N8").Value = 983.05
$newSheet.Range("O8").Value = "HOME"

$newSheet.Activate()
$newSheet.Range("B17").Select()
